# Updated TAKEN pipeline numbers
# Recalculated Total Cycles (J), Total Correct Branch Predictions (M) and
# CPI (N) for every "Taken" (Taken/Not-taken = 1) test row in the
# performance table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("J4").Value = 364569
$ws.Range("N4").Value = 10.490589999999999

$ws.Range("J6").Value = 196512
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 5.6546960000000004

$ws.Range("J8").Value = 261852
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 7.5348759999999997

$ws.Range("J10").Value = 77291
$ws.Range("M10").Value = 1277
$ws.Range("N10").Value = 2.2240730000000002

$ws.Range("J12").Value = 196512
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = 5.6546960000000004

$ws.Range("J14").Value = 77291
$ws.Range("M14").Value = 1277
$ws.Range("N14").Value = 2.2240730000000002

$ws.Range("J16").Value = 127267
$ws.Range("M16").Value = 1250
$ws.Range("N16").Value = 3.6621489999999999

$ws.Range("J18").Value = 314529
$ws.Range("M18").Value = 0
$ws.Range("N18").Value = 9.0506729999999997

$ws.Range("J20").Value = 363696
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = 10.46547

$ws.Range("J26").Value = 50133
$ws.Range("M26").Value = 1292
$ws.Range("N26").Value = 1.442593

$ws.Range("J32").Value = 47594
$ws.Range("M32").Value = 1292
$ws.Range("N32").Value = 1.3695329999999999

# Leave the final selection on the last-edited cell (matches the
# workbook's saved cursor position).
$ws.Activate() | Out-Null
$ws.Range("N32").Select() | Out-Null
